$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 246.875
$ws.Range("J18").Value = 666.6667
$ws.Range("L18").Value = 666.6667
$ws.Range("N18").Value = -1234.6667
$ws.Range("H33").Value = 178.4
$ws.Range("I33").Value = 178.4
$ws.Range("K33").Value = 178.4
$ws.Range("M33").Value = 50.59999999999999
$ws.Range("H40").Value = 1450.8125
$ws.Range("I40").Value = 1213.7273
$ws.Range("J40").Value = 1972.4
$ws.Range("K40").Value = 1213.7273
$ws.Range("L40").Value = 1972.4
$ws.Range("M40").Value = -1038.7273
$ws.Range("N40").Value = -2322.4
$ws.Range("H64").Value = 3290.3845
$ws.Range("I64").Value = 2862.5
$ws.Range("J64").Value = 3975
$ws.Range("K64").Value = 2862.5
$ws.Range("L64").Value = 3975
$ws.Range("M64").Value = -2614.5
$ws.Range("N64").Value = -4471
$ws.Range("H67").Value = 3290.3845
$ws.Range("I67").Value = 2862.5
$ws.Range("J67").Value = 3975
$ws.Range("K67").Value = 2862.5
$ws.Range("L67").Value = 3975
$ws.Range("M67").Value = -2004.5
$ws.Range("N67").Value = -5691
$ws.Range("H116").Value = 6200.125
$ws.Range("I116").Value = 2797.5
$ws.Range("J116").Value = 7334.3335
$ws.Range("K116").Value = 2797.5
$ws.Range("L116").Value = 7334.3335
$ws.Range("M116").Value = 644.5
$ws.Range("N116").Value = -14218.3335
$ws.Range("H125").Value = 505.16666
$ws.Range("I125").Value = 610.3333
$ws.Range("K125").Value = 5492.9997
$ws.Range("M125").Value = -3032.9997
$ws.Range("H138").Value = 2191.2964
$ws.Range("I138").Value = 2452.6924
$ws.Range("J138").Value = 2141.3235
$ws.Range("K138").Value = 7358.0772
$ws.Range("L138").Value = 6423.970499999999
$ws.Range("M138").Value = -2218.0772
$ws.Range("N138").Value = -16703.9705

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4522.75
$ws.Range("I32").Value = 4004.8147
$ws.Range("J32").Value = 18507
$ws.Range("K32").Value = 4004.8147
$ws.Range("L32").Value = 18507
$ws.Range("M32").Value = -3717.8147
$ws.Range("N32").Value = -19081
$ws.Range("H74").Value = 83334400
$ws.Range("I74").Value = 100000680
$ws.Range("K74").Value = 100000680
$ws.Range("M74").Value = -99999806
$ws.Range("H77").Value = 83334400
$ws.Range("I77").Value = 100000680
$ws.Range("K77").Value = 500003400
$ws.Range("M77").Value = -499999032
$ws.Range("H97").Value = 1372.625
$ws.Range("I97").Value = 1563.5
$ws.Range("J97").Value = 800
$ws.Range("K97").Value = 1563.5
$ws.Range("L97").Value = 800
$ws.Range("M97").Value = -1067.5
$ws.Range("N97").Value = -1792
$ws.Range("H110").Value = 801.75
$ws.Range("I110").Value = 692.8182
$ws.Range("K110").Value = 692.8182
$ws.Range("M110").Value = 1352.1818
$ws.Range("H132").Value = 17009.06
$ws.Range("I132").Value = 1908.3334
$ws.Range("J132").Value = 57277.668
$ws.Range("K132").Value = 5725.0002
$ws.Range("L132").Value = 171833.004
$ws.Range("M132").Value = -3195.0002
$ws.Range("N132").Value = -176893.004
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 496.66666
$ws.Range("I22").Value = 496.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 496.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -323.66666
$ws.Range("H107").Value = 2362.6
$ws.Range("I107").Value = 900
$ws.Range("J107").Value = 3337.6667
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 3337.6667
$ws.Range("M107").Value = 1020
$ws.Range("N107").Value = -7177.6667
$ws.Range("H130").Value = 55998.75
$ws.Range("J130").Value = 55998.75
$ws.Range("L130").Value = 55998.75
$ws.Range("N130").Value = -66038.75
$ws.Range("H134").Value = 5896.0586
$ws.Range("I134").Value = 6795.2144
$ws.Range("J134").Value = 1700
$ws.Range("K134").Value = 20385.6432
$ws.Range("L134").Value = 5100
$ws.Range("M134").Value = -17850.6432
$ws.Range("N134").Value = -10170

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 40020.77
$ws.Range("I58").Value = 1828.7778
$ws.Range("K58").Value = 1828.7778
$ws.Range("M58").Value = -1625.7778
$ws.Range("H62").Value = 5287.143
$ws.Range("I62").Value = 5999.6665
$ws.Range("J62").Value = 4752.75
$ws.Range("K62").Value = 5999.6665
$ws.Range("L62").Value = 4752.75
$ws.Range("M62").Value = -5375.6665
$ws.Range("N62").Value = -6000.75
$ws.Range("H65").Value = 5287.143
$ws.Range("I65").Value = 5999.6665
$ws.Range("J65").Value = 4752.75
$ws.Range("K65").Value = 29998.3325
$ws.Range("L65").Value = 23763.75
$ws.Range("M65").Value = -26878.3325
$ws.Range("N65").Value = -30003.75
$ws.Range("H105").Value = 2512.2
$ws.Range("I105").Value = 300
$ws.Range("J105").Value = 3987
$ws.Range("K105").Value = 300
$ws.Range("L105").Value = 3987
$ws.Range("M105").Value = 1447
$ws.Range("N105").Value = -7481
$ws.Range("H136").Value = 40020.77
$ws.Range("I136").Value = 1828.7778
$ws.Range("K136").Value = 5486.3334
$ws.Range("M136").Value = -2936.3334

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 168.66667
$ws.Range("I12").Value = 50
$ws.Range("J12").Value = 192.4
$ws.Range("K12").Value = 150
$ws.Range("L12").Value = 577.2
$ws.Range("M12").Value = 23
$ws.Range("N12").Value = -923.2
$ws.Range("H62").Value = 6043.143
$ws.Range("J62").Value = 9453.5
$ws.Range("L62").Value = 28360.5
$ws.Range("N62").Value = -29732.5
$ws.Range("H65").Value = 6043.143
$ws.Range("J65").Value = 9453.5
$ws.Range("L65").Value = 85081.5
$ws.Range("N65").Value = -91945.5
$ws.Range("H80").Value = 4000
$ws.Range("J80").Value = 4000
$ws.Range("L80").Value = 12000
$ws.Range("N80").Value = -13872
$ws.Range("H83").Value = 4000
$ws.Range("J83").Value = 4000
$ws.Range("L83").Value = 36000
$ws.Range("N83").Value = -45360
$ws.Range("H131").Value = 720.49
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 720.49
$ws.Range("K131").Value = 0
$ws.Range("L131").ClearContents()
$ws.Range("M131").Value = 2161.47
$ws.Range("N131").Value = -12241.47
$ws.Range("H133").Value = 5680
$ws.Range("I133").Value = 1290
$ws.Range("J133").Value = 5993.5713
$ws.Range("K133").Value = 3870
$ws.Range("L133").Value = 17980.7139
$ws.Range("M133").Value = 1190
$ws.Range("N133").Value = -28100.7139

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 29999.5
$ws.Range("J39").Value = 29999.5
$ws.Range("L39").Value = 29999.5
$ws.Range("N39").Value = -31063.5
$ws.Range("H126").Value = 5993.484
$ws.Range("I126").Value = 5479.9
$ws.Range("K126").Value = 16439.7
$ws.Range("M126").Value = -13969.7

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5138
$ws.Range("I7").Value = 3524
$ws.Range("J7").Value = 7045.4546
$ws.Range("K7").Value = 3524
$ws.Range("L7").Value = 7045.4546
$ws.Range("M7").Value = -3412
$ws.Range("N7").Value = -7269.4546
$ws.Range("H40").Value = 5279.625
$ws.Range("I40").Value = 3256.5557
$ws.Range("K40").Value = 3256.5557
$ws.Range("M40").Value = -3120.5557
$ws.Range("H126").Value = 5138
$ws.Range("I126").Value = 3524
$ws.Range("J126").Value = 7045.4546
$ws.Range("K126").Value = 10572
$ws.Range("L126").Value = 21136.3638
$ws.Range("M126").Value = -8102
$ws.Range("N126").Value = -26076.3638
$ws.Range("H136").Value = 1072.7273
$ws.Range("I136").Value = 1000
$ws.Range("J136").Value = 1133.3334
$ws.Range("K136").Value = 3000
$ws.Range("L136").Value = 3400.0002
$ws.Range("M136").Value = -450
$ws.Range("N136").Value = -8500.0002

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1164.6364
$ws.Range("I126").Value = 1055.5
$ws.Range("J126").Value = 1188.8889
$ws.Range("K126").Value = 3166.5
$ws.Range("L126").Value = 3566.6667
$ws.Range("M126").Value = -696.5
$ws.Range("N126").Value = -8506.6667
$ws.Range("H136").Value = 38463530
$ws.Range("I136").Value = 47620656
$ws.Range("K136").Value = 142861968
$ws.Range("M136").Value = -142859418
